{"js": "// 1) Drop the stale `<w:rFonts w:hint=\"cs\"/>` from the FIRST paragraph's\n//    paragraph-MARK properties (w:pPr/w:rPr) -- note the run's own rPr\n//    legitimately keeps its `w:rFonts w:hint=\"cs\"` and must stay untouched.\n//    Office.js's Font object has no \"hint\" property to toggle, so we read\n//    the paragraph's own OOXML, surgically drop just that one element from\n//    the w:pPr/w:rPr block, and write it back with insertOoxml(..., Replace)\n//    scoped to that paragraph's own range only.\nconst body = context.document.body;\nconst firstPara = body.paragraphs.getFirst();\nconst firstRange = firstPara.getRange();\n\nconst firstOoxml = firstPara.getOoxml();\nawait context.sync();\n\nconst fullPkg = firstOoxml.value;\nconst bodyMatch = fullPkg.match(/<w:body>([\\s\\S]*)<\\/w:body>/);\nconst bodyInner = bodyMatch[1];\n\n// This package always represents the target paragraph as the first <w:p>.\nconst pMatch = bodyInner.match(/<w:p\\b[\\s\\S]*?<\\/w:p>/);\nlet pXml = pMatch[0];\n\n// Drop the transient w14:paraId / w14:textId minted just for getOoxml().\npXml = pXml\n  .replace(/\\s+w14:paraId=\"[^\"]*\"/, \"\")\n  .replace(/\\s+w14:textId=\"[^\"]*\"/, \"\");\n\n// Remove the <w:rFonts .../> living inside <w:pPr><w:rPr>...</w:rPr></w:pPr>\n// (the paragraph mark's run properties) without touching the identical\n// element inside the paragraph's actual <w:r> run(s).\npXml = pXml.replace(\n  /(<w:pPr>[\\s\\S]*?<w:rPr>)([\\s\\S]*?)(<\\/w:rPr>[\\s\\S]*?<\\/w:pPr>)/,\n  (all, pre, rprInner, post) => pre + rprInner.replace(/<w:rFonts\\b[^/]*\\/>/, \"\") + post\n);\n\nconst firstParaOoxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>${pXml}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nfirstRange.insertOoxml(firstParaOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Append a new right-aligned, bidi Persian-context paragraph with the\n//    new English sentence after the last paragraph (before the sectPr).\n//    Using the body's last paragraph as the anchor makes the new run\n//    naturally inherit the same rPr (sz/szCs/lang bidi) as the paragraph\n//    it follows, matching the target markup exactly.\nconst lastPara = body.paragraphs.getLast();\nlastPara.insertParagraph(\n  \"You know what it\\u2019s very important to know the way of working with gitup.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Drop the stale <w:rFonts w:hint=\"cs\"/> from the FIRST paragraph's\n#    paragraph-mark properties (w:pPr/w:rPr). The Word object model has no\n#    \".Hint\" property to toggle, so the exact paragraph is rewritten via\n#    Range.InsertXML (a FlatOpc OOXML package) scoped to its own Range only\n#    -- InsertXML replaces just that range's contents, leaving the run's\n#    own rPr (which legitimately still carries the hint) untouched.\n$p1 = $d.Paragraphs(1)\n$firstParaXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w:rsidR=\"00743137\" w:rsidRDefault=\"0022122E\" w:rsidP=\"0022122E\"><w:pPr><w:bidi/><w:jc w:val=\"center\"/><w:rPr><w:sz w:val=\"48\"/><w:szCs w:val=\"48\"/><w:rtl/><w:lang w:bidi=\"fa-IR\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=\"cs\"/><w:sz w:val=\"48\"/><w:szCs w:val=\"48\"/><w:rtl/><w:lang w:bidi=\"fa-IR\"/></w:rPr><w:t>\u0628\u0633\u0645 \u0627\u0644\u0644\u0647 \u0627\u0644\u0631\u062d\u0645\u0646 \u0627\u0644\u0631\u062d\u06cc\u0645</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n$p1.Range.InsertXML($firstParaXml)\n\n# 2) Append a new right-aligned, bidi-Persian-context paragraph with the new\n#    English sentence at the very end of the document body (after the last\n#    existing paragraph, before the section break). Collapsing the end of\n#    the document's Content range and inserting there makes the new\n#    paragraph/run naturally inherit the same rPr (sz/szCs/lang bidi) as\n#    the paragraph it follows, matching the target markup exactly.\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n$end.Collapse(0)\n$end.InsertAfter(\"You know what it\u2019s very important to know the way of working with gitup.\")\n"}
